$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 1.1770765782808947
$ws.Range("C2").Value = 1.3708319777246605
$ws.Range("D2").Value = 3.4592388228240232
$ws.Range("E2").Value = 1.0266964026154932

$ws.Range("B3").Value = 1.8997263969175724
$ws.Range("C3").Value = 0.79165958231456479
$ws.Range("D3").Value = 3.2769210066465044
$ws.Range("E3").Value = 0.27686809088726183

$ws.Range("B1:E3").Select()
